# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("G2").Value  = -0.01948648135249794
$ws.Range("H2").Value  = -0.01948648135249794
$ws.Range("I2").Value  = -0.01206230608497436
$ws.Range("J2").Value  = -0.01206230608497436
$ws.Range("K2").Value  = 0.742
$ws.Range("L2").Value  = 0.0001174571012473881
$ws.Range("M2").Value  = 0
$ws.Range("N2").Value  = 0
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 0
$ws.Range("Q2").Value  = 0
$ws.Range("R2").Value  = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value  = 2334.2
$ws.Range("V2").Value  = 0.6989878421273282
$ws.Range("W2").Value  = 0.0001047549130336571
$ws.Range("X2").Value  = 0.09509188741269362
$ws.Range("Y2").Value  = -0.09498713249965995
$ws.Range("Z2").Value  = 1.225612689937761
$ws.Range("AA2").Value = -0.01478371540765804
$ws.Range("AB2").Value = 0.09509188741269362
$ws.Range("AC2").Value = -0.1098756028203517
$ws.Range("AG2").Value = -2334.2
$ws.Range("AJ2").Value = -2.322124950258654
$ws.Range("AK2").Value = -0.657650804383963
$ws.Range("AL2").Value = 3.73
$ws.Range("AM2").Value = 3.73
$ws.Range("AN2").Value = -0.0
$ws.Range("AO2").Value = -20.42895442359249
$ws.Range("AP2").Value = 31.08122503328895
$ws.Range("AQ2").Value = -20.42895442359249

# ---- Row 3 ----
$ws.Range("G3").Value  = -0.01948648135249794
$ws.Range("H3").Value  = -0.01948648135249794
$ws.Range("I3").Value  = -0.01206230608497436
$ws.Range("J3").Value  = -0.01206230608497436
$ws.Range("K3").Value  = 0.742
$ws.Range("L3").Value  = 0.0001174571012473881
$ws.Range("M3").Value  = -0.0
$ws.Range("N3").Value  = -0.0
$ws.Range("O3").Value  = -0.0
$ws.Range("P3").Value  = -0.0
$ws.Range("Q3").Value  = -0.0
$ws.Range("R3").Value  = -0.0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value  = 2334.2
$ws.Range("V3").Value  = 0.6989878421273282
$ws.Range("W3").Value  = 0.0001047549130336571
$ws.Range("X3").Value  = 0.09509188741269362
$ws.Range("Y3").Value  = -0.09498713249965995
$ws.Range("Z3").Value  = 1.225612689937761
$ws.Range("AA3").Value = -0.01478371540765804
$ws.Range("AB3").Value = 0.09509188741269362
$ws.Range("AC3").Value = -0.1098756028203517
$ws.Range("AG3").Value = -2334.2
$ws.Range("AJ3").Value = -2.322124950258654
$ws.Range("AK3").Value = -0.657650804383963
$ws.Range("AL3").Value = 3.73
$ws.Range("AM3").Value = 3.73
$ws.Range("AN3").Value = -0.0
$ws.Range("AO3").Value = -20.42895442359249
$ws.Range("AP3").Value = 31.08122503328895
$ws.Range("AQ3").Value = -20.42895442359249
